# Update occupation classification code:
#  - insert a new "% of total people" column between "# of people" and
#    "average debt per person"
#  - populate it with each occupation's share of the total people count
#  - re-sort the occupation rows descending by 6p_total (column C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column E; old E (average debt per person) shifts to F.
$ws.Columns.Item(5).Insert()
$ws.Range("E1").Value = "% of total people"

# 2) Compute the % of total people for each occupation row (rows 2-18).
$total = 0
for ($r = 2; $r -le 18; $r++) {
    $total = $total + $ws.Cells.Item($r, 4).Value()
}
for ($r = 2; $r -le 18; $r++) {
    $count = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 5).Value = ($count / $total) * 100
}

# 3) Re-sort the data rows (B2:F18) descending by 6p_total (column C),
#    leaving the running index in column A untouched.
$sortRange = $ws.Range("B2:F18")
$sortKey = $ws.Range("C2")
$sortRange.Sort($sortKey, 2)
